$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9485791610284168
$ws.Range("C2").Value = 0.9576502732240437
$ws.Range("D2").Value = 0.9530931339225017

$ws.Range("B3").Value = 0.9592641261498029
$ws.Range("C3").Value = 0.9505208333333334
$ws.Range("D3").Value = 0.9548724656638327

$ws.Range("B5").Value = 0.9539216435891098
$ws.Range("C5").Value = 0.9540855532786885
$ws.Range("D5").Value = 0.9539827997931671

$ws.Range("B6").Value = 0.9540498631705665
$ws.Range("D6").Value = 0.9540041517740632
